# Ajouts images déprt CLMI Cité de l'Or
# Swap the "Sécurité Valcourt" / "Valcourt Securtiy services" row (row 16)
# on the ADMIN sheet so column A holds the service description and column B
# holds the (wrapped) label, matching the rest of the table's layout, and
# move the active selection to A21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADMIN")
$ws.Activate() | Out-Null

$a16 = $ws.Range("A16")
$b16 = $ws.Range("B16")

$aVal = $a16.Value2
$bVal = $b16.Value2

$a16.Value = $bVal
$b16.Value = $aVal

$a16.WrapText = $true
$b16.WrapText = $false

$ws.Range("A21").Select() | Out-Null
